$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new log entry row (row 59) below the existing data (which ends at row 58).
# Column A: description of the new work item (becomes a new shared string).
# Column B: count of items/papers.
# Column C: time spent (minutes).
$ws.Cells.Item(59, 1).Value = "Create new figures for RQ2.1"
$ws.Cells.Item(59, 2).Value = 1
$ws.Cells.Item(59, 3).Value = 20

# Update the on-screen selection to the newly added cells.
$ws.Range("C57:C59").Select() | Out-Null
